# Outstandings.xlsx edit script
# Summary of change: on "Purchase 22-23" sheet, a new invoice line (25/23-24,
# Namrata Rubber Product Pvt Ltd, dated 45099, amount 40179) was inserted as a
# new row 15 inside the "Sr. No 1" group (right after the 45093 / 23/23-24
# line). All subsequent rows shift down by one. The amount on the existing
# 45093 line (row 14) was also corrected from 87556.5 to 87556 (whole number).
# Minor view/column-width tweaks on both sheets are also reproduced.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Purchase 22-23" --------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Correct the amount that was previously 87556.5 -> 87556 (row 14, col E)
$ws1.Range("E14").Value = 87556

# Insert a brand-new row at position 15; this shifts every row >= 15 down by
# one and Excel automatically rewrites the dependent SUM-style formulas
# (F21 -> F22 etc.) to point at the shifted cells.
$ws1.Rows("15:15").Insert()

# The freshly inserted row 15 inherits blank/default formatting; pull in the
# same per-column look used by the rest of this data block (row 17, which
# used to be row 16 before the insert) so number formats / borders / fonts
# line up with the surrounding rows.
$ws1.Rows("17:17").Copy()
$ws1.Rows("15:15").PasteSpecial(-4122)
$ws1.Range("G15:XFD15").Clear()

# Populate the new invoice line.
$ws1.Range("A15").Value = ""
$ws1.Range("B15").Value = 45099
$ws1.Range("C15").Value = "25/23-24"
$ws1.Range("D15").Value = "Namrata Rubber Product Pvt Ltd"
$ws1.Range("E15").Value = 40179
$ws1.Range("F15").Formula = "=F14+E15"

# F15 is now the running-total cell that closes out the "Sr. No 1" block, so
# it should carry the bold "group total" look (same as F14 used to have).
$ws1.Range("F15").Font.Bold = $true
$ws1.Range("F15").Font.Name = "Arial"
$ws1.Range("F15").Font.Size = 10

# F14 is no longer the last line of its group, so it reverts to the plain
# (non-bold) running-total style used by the rest of the block (matches F13).
$ws1.Range("F13").Copy()
$ws1.Range("F14").PasteSpecial(-4122)
$ws1.Range("F14").Formula = "=F13+E14"

# Restore the view: selection moves to A16, and the previous "scrolled down"
# top-left cell is no longer pinned.
$ws1.Range("A16").Select()

# --- Sheet 2 ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Narrower first three columns.
$ws2.Columns.Item(1).ColumnWidth = 7.5
$ws2.Columns.Item(2).ColumnWidth = 13
$ws2.Columns.Item(3).ColumnWidth = 16.83
